# Apply updated robustness-test rows to the active worksheet.
# NOTE: column E (error_type) values are unchanged by this edit for every
# row, so they are intentionally left untouched below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, message_id (A), message (B), sentence (C), boundary (D)
$rows = @(
    @(2,  3,  "Low Power Returning home .", "Low Power", "0-1"),
    @(3,  3,  "Low Power Returning home .", "Returning home", "2-3"),
    @(4,  3,  "Low Power Returning home .", "Low Power Returning home", "0-3"),
    @(5,  6,  "Aircraft ActiveTrack available at max speed When exceeding nnn, Obstacle Avoidance is not available .", "When exceeding nnn, Obstacle Avoidance is not available", "6-13"),
    @(6,  6,  "Aircraft ActiveTrack available at max speed When exceeding nnn, Obstacle Avoidance is not available .", "When exceeding nnn,", "6-8"),
    @(7,  13, "Aircraft is close to the Home Point Initiating Return to Home will now trigger Auto Landing .", "Initiating Return to Home will now trigger Auto Landing", "7-15"),
    @(8,  13, "Aircraft is close to the Home Point Initiating Return to Home will now trigger Auto Landing .", "Return to Home will now trigger Auto Landing", "8-15"),
    @(9,  15, "Critically Low Voltage Aircraft will land .", "Critically Low Voltage", "0-2"),
    @(10, 15, "Critically Low Voltage Aircraft will land .", "Aircraft will land", "3-5"),
    @(11, 15, "Critically Low Voltage Aircraft will land .", "Critically Low Voltage Aircraft will land", "0-5"),
    @(12, 21, "Propeller Guard Mounted Propeller Guard mounted Forward Obstacle Sensing will be automatically switched off .", "Propeller Guard Mounted Propeller Guard mounted", "0-5"),
    @(13, 21, "Propeller Guard Mounted Propeller Guard mounted Forward Obstacle Sensing will be automatically switched off .", "Propeller Guard Mounted", "0-2"),
    @(14, 21, "Propeller Guard Mounted Propeller Guard mounted Forward Obstacle Sensing will be automatically switched off .", "Propeller Guard mounted", "3-5")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
